$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-260). The commit bumps that date by one day (45181 -> 45182)
# across all of them.
$ws.Range("C2:C260").Value = 45182
